$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.309.89"
$ws.Range("E2").Value = "  -2.71%  "

$ws.Range("D3").Value = "1.570.29"
$ws.Range("E3").Value = "  -3.72%  "

$ws.Range("E4").Value = "  -0.40%  "

$ws.Range("D5").Value = "'207.91"
$ws.Range("E5").Value = "  -3.03%  "

$ws.Range("E6").Value = "  -0.37%  "

$ws.Range("E7").Value = "  -5.02%  "

$ws.Range("E8").Value = "  -2.16%  "

$ws.Range("D9").Value = "'0.0607"
$ws.Range("E9").Value = "  -1.70%  "

$ws.Range("D10").Value = "'17.95"
$ws.Range("E10").Value = "  -2.02%  "

$ws.Range("D11").Value = "'0.0783"
$ws.Range("E11").Value = "  -0.94%  "

$ws.Range("D12").Value = "1.788.89"
$ws.Range("E12").Value = "  -3.75%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.03"
$ws.Range("E13").Value = "  -3.01%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.562.04"
$ws.Range("E14").Value = "  -4.25%  "

$ws.Range("E15").Value = "  -3.05%  "

$ws.Range("D16").Value = "25.303.00"
$ws.Range("E16").Value = "  -2.64%  "

$ws.Range("D17").Value = "'59.74"
$ws.Range("E17").Value = "  -2.49%  "

$ws.Range("E18").Value = "  -4.20%  "

$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").Value = "'185.86"
$ws.Range("E20").Value = "  -2.08%  "

$ws.Range("E21").Value = "  -2.23%  "

$ws.Range("D22").Value = "'9.33"
$ws.Range("E22").Value = "  -2.52%  "

$ws.Range("D23").Value = "'5.89"
$ws.Range("E23").Value = "  -2.93%  "

$ws.Range("E24").Value = "  -2.46%  "

$ws.Range("E25").Value = "  -0.33%  "

$ws.Range("D26").Value = "'141.01"
$ws.Range("E26").Value = "  -1.75%  "

$ws.Range("D27").Value = "'1.67"
$ws.Range("E27").Value = "  -6.58%  "

$ws.Range("D28").Value = "'6.46"
$ws.Range("E28").Value = "  -3.66%  "

$ws.Range("E29").Value = "  -1.67%  "

$ws.Range("E30").Value = "  -6.27%  "

$ws.Range("D31").Value = "'0.0463"
$ws.Range("E31").Value = "  -3.75%  "

$ws.Range("D32").Value = "'3.06"

$ws.Range("E33").Value = "  -3.45%  "

$ws.Range("E34").Value = "  -1.16%  "

$ws.Range("D35").Value = "'2.30"
$ws.Range("E35").Value = "  -4.54%  "

$ws.Range("D36").Value = "1.095.40"
$ws.Range("E36").Value = "  -3.28%  "

$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("E38").Value = "  -2.31%  "

$ws.Range("D39").Value = "'2.32"
$ws.Range("E39").Value = "  -4.98%  "

$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'0.496"
$ws.Range("E40").Value = "  -3.70%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'0.776"
$ws.Range("E41").Value = "  -8.81%  "

$ws.Range("D42").Value = "'0.770"
$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("D43").Value = "'92.62"
$ws.Range("E43").Value = "  -5.52%  "

$ws.Range("E44").Value = "  -2.48%  "

$ws.Range("D45").Value = "1.703.24"
$ws.Range("E45").Value = "  -3.72%  "

$ws.Range("E46").Value = "  -2.76%  "

$ws.Range("D47").Value = "'52.91"
$ws.Range("E47").Value = "  -3.16%  "

$ws.Range("D48").Value = "'0.0506"
$ws.Range("E48").Value = "  -3.63%  "

$ws.Range("E49").Value = "  -3.60%  "

$ws.Range("D50").Value = "'0.406"
$ws.Range("E50").Value = "  -2.10%  "

$ws.Range("E51").Value = "  -0.51%  "
